$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cells value while preserving its original style/format.
# Plain numeric-looking strings (e.g. "336.98") would otherwise be auto-
# converted to numbers by Excel; briefly forcing a text NumberFormat keeps
# them as text (matching the source inlineStr cells), then the cells
# original style is restored so no visible formatting changes.
function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

$ws.Range('D2').Value = '28.037.77'
$ws.Range('E2').Value = '  -0.20%  '
$ws.Range('D3').Value = '1.815.40'
$ws.Range('E3').Value = '  +2.26%  '
$ws.Range('E4').Value = '  -0.38%  '
Set-TextValue 'D5' '336.98'
$ws.Range('E5').Value = '  -0.61%  '
Set-TextValue 'D6' '1.000'
$ws.Range('E6').Value = '  -0.55%  '
Set-TextValue 'D7' '0.4203'
$ws.Range('E7').Value = '  +9.76%  '
Set-TextValue 'D8' '0.3522'
$ws.Range('E8').Value = '  +2.87%  '
Set-TextValue 'D9' '45.48'
$ws.Range('E9').Value = '  -3.40%  '
Set-TextValue 'D10' '1.158'
$ws.Range('E10').Value = '  +0.67%  '
Set-TextValue 'D11' '0.07542'
$ws.Range('E11').Value = '  +1.90%  '
Set-TextValue 'D12' '22.84'
$ws.Range('E12').Value = '  -3.26%  '
Set-TextValue 'D13' '1.002'
$ws.Range('E13').Value = '  -0.44%  '
Set-TextValue 'D14' '6.295'
$ws.Range('E14').Value = '  -2.13%  '
Set-TextValue 'D15' '7.280'
$ws.Range('E15').Value = '  -0.72%  '
$ws.Range('D16').Value = '1.813.48'
$ws.Range('E16').Value = '  +1.64%  '
Set-TextValue 'D17' '0.00001091'
$ws.Range('E17').Value = '  +1.20%  '
Set-TextValue 'D18' '0.06687'
$ws.Range('E18').Value = '  -0.08%  '
Set-TextValue 'D19' '82.59'
$ws.Range('E19').Value = '  +0.24%  '
Set-TextValue 'D20' '1.001'
$ws.Range('E20').Value = '  -0.35%  '
Set-TextValue 'D21' '17.40'
$ws.Range('E21').Value = '  -0.13%  '
Set-TextValue 'D22' '6.387'
$ws.Range('E22').Value = '  -0.43%  '
$ws.Range('D23').Value = '28.097.29'
$ws.Range('E23').Value = '  -0.01%  '
$ws.Range('E24').Value = '  -1.54%  '
Set-TextValue 'D25' '2.396'
$ws.Range('E25').Value = '  +0.03%  '
Set-TextValue 'D26' '2.506'
$ws.Range('E26').Value = '  +3.93%  '
$ws.Range('E27').Value = '  +0.42%  '
Set-TextValue 'D28' '156.31'
$ws.Range('E28').Value = '  +1.58%  '
$ws.Range('D29').Value = '2.020.86'
$ws.Range('E29').Value = '  +1.90%  '
Set-TextValue 'D30' '1.312'
$ws.Range('E30').Value = '  -7.36%  '
Set-TextValue 'D31' '133.32'
Set-TextValue 'D32' '4.080'
$ws.Range('E32').Value = '  +1.14%  '
Set-TextValue 'D33' '6.008'
$ws.Range('E33').Value = '  -1.29%  '
Set-TextValue 'D34' '0.09155'
$ws.Range('E34').Value = '  +2.81%  '
Set-TextValue 'D35' '12.38'
$ws.Range('E35').Value = '  -2.96%  '
Set-TextValue 'D36' '0.06353'
$ws.Range('E36').Value = '  -0.07%  '
Set-TextValue 'D37' '0.02352'
$ws.Range('E37').Value = '  -2.01%  '
Set-TextValue 'D38' '0.6685'
$ws.Range('E38').Value = '  -1.94%  '
Set-TextValue 'D39' '5.247'
$ws.Range('E39').Value = '  -1.40%  '
Set-TextValue 'D40' '0.2164'
$ws.Range('E40').Value = '  +0.36%  '
Set-TextValue 'D41' '1.511'
$ws.Range('E41').Value = '  +0.50%  '
Set-TextValue 'D42' '1.219'
$ws.Range('E42').Value = '  -2.24%  '
Set-TextValue 'D43' '8.158'
$ws.Range('E43').Value = '  -1.55%  '
Set-TextValue 'D44' '14.17'
$ws.Range('E44').Value = '  +0.22%  '
Set-TextValue 'D45' '1.000'
$ws.Range('E45').Value = '  -0.43%  '
Set-TextValue 'D46' '0.6166'
$ws.Range('E46').Value = '  -1.46%  '
Set-TextValue 'D47' '3.877'
$ws.Range('E47').Value = '  +0.38%  '
Set-TextValue 'D48' '128.51'
$ws.Range('E48').Value = '  -3.16%  '
Set-TextValue 'D49' '2.060'
$ws.Range('E49').Value = '  -0.19%  '
Set-TextValue 'D50' '1.186'
$ws.Range('E50').Value = '  -0.49%  '
Set-TextValue 'D51' '0.07126'
$ws.Range('E51').Value = '  -5.24%  '
